# Updated symbol list on Thu Dec 22 15:40:43 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (D), "Volume(1h)" (E) ranking labels, and (for the
# rows whose coin fell to/rose from a neighbouring rank) the "Coin" (B) /
# "Link" (C) columns on Sheet1 of the crypto snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row=2; D='241.38' }
    @{ Row=3; D='21.81' }
    @{ Row=4; D='5.357' }
    @{ Row=5; D='0.05699' }
    @{ Row=6; D='3.427' }
    @{ Row=7; D='6.292' }
    @{ Row=8; D='0.8060' }
    @{ Row=9; D='0.8510' }
    @{ Row=11; D='0.07278' }
    @{ Row=12; D='0.03029' }
    @{ Row=13; D='0.03151' }
    @{ Row=14; B='ProBitToken'; C='https://coinranking.com/coin/lQP4d6T2+probittoken-prob'; D='0.1242'; E='13ProBitTokenPROB' }
    @{ Row=15; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.09361'; E='14BitMartTokenBMX' }
    @{ Row=16; B='MCDex'; C='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'; D='3.926'; E='15MCDexMCB' }
    @{ Row=17; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001586'; E='16BitForexTokenBF' }
    @{ Row=18; B='CoinExToken'; C='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'; D='0.04828'; E='17CoinExTokenCET' }
    @{ Row=19; D='0.006322' }
    @{ Row=20; D='0.0009970' }
    @{ Row=21; D='0.004055' }
    @{ Row=23; D='3.717' }
    @{ Row=24; D='2.166' }
    @{ Row=25; B='One'; C='https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; D='0.01073'; E='24OneONEBestin24h' }
    @{ Row=26; B='BitpandaEcosystemToken'; C='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'; D='0.3232'; E='25BitpandaEcosystemTokenBEST' }
    @{ Row=27; D='0.0003032' }
    @{ Row=40; D='0.03820' }
    @{ Row=41; D='0.006759' }
    @{ Row=43; D='0.003206' }
    @{ Row=44; D='0.006366' }
    @{ Row=45; D='0.00005618' }
    @{ Row=47; D='0.5810'; E='46CoinbaseStockTokenCOIN' }
    @{ Row=48; D='0.1423' }
    @{ Row=49; D='0.00002104' }
    @{ Row=50; D='0.01012' }
)

foreach ($row in $changes) {
    if ($row.ContainsKey('B')) {
        $ws.Cells.Item($row.Row, 2).Value = $row.B
    }
    if ($row.ContainsKey('C')) {
        $ws.Cells.Item($row.Row, 3).Value = $row.C
    }
    if ($row.ContainsKey('D')) {
        # Column D holds the price as plain text in this workbook (the source
        # feed writes it as an inline string, not a number) - force the Text
        # number format before assigning so the value is stored as text
        # instead of Excel auto-coercing the numeric-looking string into a
        # real number.
        $cell = $ws.Cells.Item($row.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $row.D
    }
    if ($row.ContainsKey('E')) {
        $ws.Cells.Item($row.Row, 5).Value = $row.E
    }
}

Write-Host "Updated $($changes.Count) rows"
